# refs #166 * Software Architektur
#
# Applies the architecture-overview slide changes:
#  - merge the "Android" / "-Client" runs into a single "Android-Client" run
#  - remove the "Aussendienstmitarbeiter" textbox and its dashed connector
#  - rename/resize the three "Application" round-rect boxes
#    ("Ruby on Rails Application", "Android Application", "Browser")
#  - re-route the two arrow connectors attached to those boxes

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) "Android" + "-Client" -> single run "Android-Client" (drop err flag)
# ---------------------------------------------------------------------
$androidClient = $s.Shapes.Item("Textfeld 21")
$tr = $androidClient.TextFrame.TextRange
$prefix = $tr.Characters(1, 7)        # "Android"
$prefix.Text = ""                      # leaves the clean "-Client" run
$tr.InsertBefore("Android")            # re-add prefix using clean formatting
# force a real re-merge of the (now identically formatted) runs; use a
# placeholder with no shared prefix/substring so the engine rebuilds the
# run list instead of doing a minimal partial-text patch
$tr.Text = "zzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$tr.Text = "Android-Client"

# ---------------------------------------------------------------------
# 2) Remove the "Aussendienstmitarbeiter" textbox ...
# ---------------------------------------------------------------------
$s.Shapes.Item("Textfeld 22").Delete()

# ...and the dashed connector that used to link it to the Android textbox
$s.Shapes.Item("Gerade Verbindung 1060").Delete()

# ---------------------------------------------------------------------
# 3) "Application" -> "Ruby on Rails Application" round rect, resize
# ---------------------------------------------------------------------
$rails = $s.Shapes.Item("Abgerundetes Rechteck 1036")
$rails.Left = 251.40866861732283
$rails.Top = 191.77889763779527
$rails.Width = 128.98125984251968
$rails.Height = 44.201574803149605
$rails.TextFrame.TextRange.Text = "Ruby on Rails Application"

# ---------------------------------------------------------------------
# 4) "Application" -> "Android Application" round rect, resize
# ---------------------------------------------------------------------
$androidApp = $s.Shapes.Item("Abgerundetes Rechteck 45")
$androidApp.Left = 342.92826851653547
$androidApp.Top = 304.0195275590551
$androidApp.Width = 88.50220472440945
$androidApp.Height = 34.798976377952755
$androidApp.TextFrame.TextRange.Text = "Android Application"

# ---------------------------------------------------------------------
# 5) Connector between "Android Application" and "Ruby on Rails
#    Application" follows the resized boxes
# ---------------------------------------------------------------------
$conn1059 = $s.Shapes.Item("Gerade Verbindung mit Pfeil 1058")
$conn1059.Left = 315.89929133858266
$conn1059.Top = 235.98047644094487
$conn1059.Width = 71.28007894015748
$conn1059.Height = 68.03905871811024

# ---------------------------------------------------------------------
# 6) "Application" -> "Browser" round rect, resize
# ---------------------------------------------------------------------
$browser = $s.Shapes.Item("Abgerundetes Rechteck 112")
$browser.Left = 185.28952795905514
$browser.Top = 304.0195275590551
$browser.Width = 88.50220472440945
$browser.Height = 35.1000004
$browser.TextFrame.TextRange.Text = "Browser"

# ---------------------------------------------------------------------
# 7) Connector between "Browser" and "Ruby on Rails Application"
#    follows the resized boxes
# ---------------------------------------------------------------------
$conn107 = $s.Shapes.Item("Gerade Verbindung mit Pfeil 106")
$conn107.Left = 229.54062992125984
$conn107.Top = 235.98047644094487
$conn107.Width = 86.35866171732283
$conn107.Height = 68.03905871811024
